$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently sits right under
#    the title heading.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# 2. Insert a new bold paragraph reading
#    "Play Betti the Yetti Slot Game Free - Review 2021" right before the
#    final paragraph (the one holding the AI-image-prompt text).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($count)
$finalPara.Range.InsertParagraphBefore() | Out-Null

$newPara = $d.Paragraphs.Item($count)
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Betti the Yetti Slot Game Free - Review 2021</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newParaXml) | Out-Null

# ---------------------------------------------------------------------------
# 3. Replace the old "Please create a cartoon-style image..." AI-art prompt
#    with the meta-description text (the run keeps its italic formatting).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Please create a cartoon-style image featuring a happy Maya warrior with glasses, fitting the game " + [char]34 + "Betti the Yetti" + [char]34 + ". The Maya warrior should be smiling, with a friendly and welcoming expression. The warrior's glasses should be visible and slightly oversized, adding to the cartoon style of the image. In the background, the Himalayan mountain should be visible, with trees framing the image. The overall style should be whimsical and fun, capturing the adventurous spirit of the game and the unique character of the Maya warrior.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Join the adventure and play Betti the Yetti slot game for free. Check our 2021 review, the bonus rounds, RTP, graphics, and sounds.",
    2) | Out-Null
